# Add 2020-07-28 data update: refresh nombre_aides (C) and montant_total (D)
# values for a set of rows in Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 6;  C = "21";   D = "124500.00" },
    @{ Row = 22; C = "305";  D = "844773.99" },
    @{ Row = 23; C = "106";  D = "383337.00" },
    @{ Row = 24; C = "34";   D = "150233.00" },
    @{ Row = 32; C = "80";   D = "190320.00" },
    @{ Row = 33; C = "442";  D = "1168397.87" },
    @{ Row = 34; C = "175";  D = "723340.11" },
    @{ Row = 35; C = "61";   D = "288974.00" },
    @{ Row = 36; C = "20";   D = "110500.00" },
    @{ Row = 44; C = "15";   D = "61121.84" },
    @{ Row = 45; C = "48";   D = "211622.07" },
    @{ Row = 46; C = "24";   D = "141780.00" },
    @{ Row = 49; C = "81";   D = "227937.17" },
    @{ Row = 50; C = "498";  D = "1525232.54" },
    @{ Row = 51; C = "219";  D = "841064.15" },
    @{ Row = 73; C = "203";  D = "498326.09" },
    @{ Row = 74; C = "799";  D = "2306575.70" },
    @{ Row = 76; C = "96";   D = "397484.52" },
    @{ Row = 77; C = "23";   D = "130383.20" },
    @{ Row = 78; C = "23";   D = "47500.00" },
    @{ Row = 85; C = "86";   D = "208800.00" },
    @{ Row = 86; C = "379";  D = "1064760.67" },
    @{ Row = 87; C = "162";  D = "610292.91" },
    @{ Row = 90; C = "9";    D = "18000.00" },
    @{ Row = 92; C = "1089"; D = "3005806.10" },
    @{ Row = 93; C = "400";  D = "1409131.02" }
)

# Prefix with an apostrophe so Excel stores these as text (matching the
# source data, which keeps numeric-looking values as text strings) rather
# than coercing them into numeric cells and dropping formatting such as
# trailing ".00".
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = "'" + $u.C
    $ws.Cells.Item($u.Row, 4).Value = "'" + $u.D
}
